# Update row 6 ("signal segment 5") data across the Step1/Step2/Step3 sheets
# to reflect the recomputed/updated mounted pipeline values.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Step1_Data")
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.2503307642664209
$ws.Range("F6").Value = 0.08608628482077046
$ws.Range("G6").Value = 0.4213607802461588
$ws.Range("H6").Value = 0.01898255824859036
$ws.Range("I6").Value = 0.02820188454478797
$ws.Range("J6").Value = 0.02152339116282657
$ws.Range("K6").Value = 0.003420034168001239
$ws.Range("O6").Value = 0.003117099404869358
$ws.Range("P6").Value = 0.0376200389145691
$ws.Range("Q6").Value = 0.01128204778021782
$ws.Range("S6").Value = 0.007837840720066248
$ws.Range("T6").Value = 0.002575266190961002
$ws.Range("U6").Value = 0.01792956367937711
$ws.Range("V6").Value = 0.001408821707461609
$ws.Range("Z6").Value = 0.00187052651833792
$ws.Range("AA6").Value = 0.006940507104306092
$ws.Range("AC6").Value = 0.01867843488307695
$ws.Range("AE6").Value = 0.01470130323844228
$ws.Range("AF6").Value = 0.02068781058239652
$ws.Range("AH6").Value = 0.01620473588585596
$ws.Range("AI6").Value = 0.005132423286430712
$ws.Range("AJ6").Value = 0.004107882646074946

$ws = $wb.Worksheets.Item("Step2_Sj")
$ws.Range("D6").Value = 0
$ws.Range("E6").Value = 0.2503307642664209
$ws.Range("F6").Value = 0.3364170490871914
$ws.Range("G6").Value = 0.7577778293333501
$ws.Range("H6").Value = 0.7767603875819405
$ws.Range("I6").Value = 0.8049622721267286
$ws.Range("J6").Value = 0.8264856632895551
$ws.Range("K6").Value = 0.8299056974575563
$ws.Range("L6").Value = 0.8299056974575563
$ws.Range("M6").Value = 0.8299056974575563
$ws.Range("N6").Value = 0.8299056974575563
$ws.Range("O6").Value = 0.8330227968624256
$ws.Range("P6").Value = 0.8706428357769948
$ws.Range("Q6").Value = 0.8819248835572125
$ws.Range("R6").Value = 0.8819248835572125
$ws.Range("S6").Value = 0.8897627242772788
$ws.Range("T6").Value = 0.8923379904682398
$ws.Range("U6").Value = 0.910267554147617
$ws.Range("V6").Value = 0.9116763758550785
$ws.Range("W6").Value = 0.9116763758550785
$ws.Range("X6").Value = 0.9116763758550785
$ws.Range("Y6").Value = 0.9116763758550785
$ws.Range("Z6").Value = 0.9135469023734165
$ws.Range("AA6").Value = 0.9204874094777226
$ws.Range("AB6").Value = 0.9204874094777226
$ws.Range("AC6").Value = 0.9391658443607995
$ws.Range("AD6").Value = 0.9391658443607995
$ws.Range("AE6").Value = 0.9538671475992418
$ws.Range("AF6").Value = 0.9745549581816384
$ws.Range("AG6").Value = 0.9745549581816384
$ws.Range("AH6").Value = 0.9907596940674943
$ws.Range("AI6").Value = 0.995892117353925

$ws = $wb.Worksheets.Item("Step3_DataPts_0.5")
$ws.Range("F6").Value = 0.7577778293333501

$ws = $wb.Worksheets.Item("Step3_DataPts_0.7")
$ws.Range("F6").Value = 0.7577778293333501

$ws = $wb.Worksheets.Item("Step3_DataPts_0.8")
$ws.Range("F6").Value = 0.8049622721267286

$ws = $wb.Worksheets.Item("Step3_DataPts_0.9")
$ws.Range("F6").Value = 0.910267554147617
